$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, matching style of other headers
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       RandomForestRegressor())]),`n                                            param_grid={'model__max_depth': [3,`n                                                                             5,`n                                                                             7],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# Update B2:D2 values
$ws.Range("B2").Value = 0.05388779317894196
$ws.Range("C2").Value = 0.9984152519378109
$ws.Range("D2").Value = 0.1693000680485397
$ws.Range("F2").Value = $modelText

# Update B3:D3 values
$ws.Range("B3").Value = 0.06140542632032733
$ws.Range("C3").Value = 0.9994175420126326
$ws.Range("D3").Value = 0.1795634552176902
$ws.Range("F3").Value = $modelText

# Update B4:D4 values
$ws.Range("B4").Value = 0.05597081437957018
$ws.Range("C4").Value = 0.9992421897229883
$ws.Range("D4").Value = 0.1918825461657122
$ws.Range("F4").Value = $modelText
